$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 1000, shifting all existing data (rows 1000-1126)
# down to rows 1002-1128. This makes room for a new weekly price report.
$ws.Range("A1000:A1001").EntireRow.Insert()

# New row 1000: "Primera" quality record for the new date
$ws.Range("A1000").Value = 6
$ws.Range("B1000").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1000").Value = "Metropolitana"
$ws.Range("D1000").Value = 45124
$ws.Range("E1000").Value = 13
$ws.Range("F1000").Value = 100112017
$ws.Range("G1000").Value = "Apio"
$ws.Range("H1000").Value = "Americana (o)"
$ws.Range("I1000").Value = "Primera"
$ws.Range("J1000").Value = 2600
$ws.Range("K1000").Value = 6000
$ws.Range("L1000").Value = 7000
$ws.Range("M1000").Value = 6462
$ws.Range("N1000").Value = "`$/docena de matas"
$ws.Range("O1000").Value = "Región de Coquimbo"
$ws.Range("P1000").Value = 1077
$ws.Range("Q1000").Value = 6
$ws.Range("R1000").Value = "Hortaliza"

# New row 1001: "Segunda" quality record for the new date
$ws.Range("A1001").Value = 6
$ws.Range("B1001").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1001").Value = "Metropolitana"
$ws.Range("D1001").Value = 45124
$ws.Range("E1001").Value = 13
$ws.Range("F1001").Value = 100112017
$ws.Range("G1001").Value = "Apio"
$ws.Range("H1001").Value = "Americana (o)"
$ws.Range("I1001").Value = "Segunda"
$ws.Range("J1001").Value = 1400
$ws.Range("K1001").Value = 4000
$ws.Range("L1001").Value = 5000
$ws.Range("M1001").Value = 4571
$ws.Range("N1001").Value = "`$/docena de matas"
$ws.Range("O1001").Value = "Región de Coquimbo"
$ws.Range("P1001").Value = 762
$ws.Range("Q1001").Value = 6
$ws.Range("R1001").Value = "Hortaliza"
